# Update crypto price/volume data to latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''25.828.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -0.03%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.642.99'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  +0.18%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  -0.26%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''215.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -0.02%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = '''  -0.53%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  -0.11%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  +0.24%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '''  -0.77%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''19.55'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -4.41%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  +0.45%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''1.656.52'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  +0.94%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''4.26'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -0.41%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''1.875.27'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +0.47%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = '''  -1.15%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''0.0₃0769'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  +0.53%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = '''  -0.54%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''25.876.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +0.09%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = '''  -0.25%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = '''  +1.44%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''194.51'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +0.99%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''9.99'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  +0.64%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  +2.07%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = '''  +0.02%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''1.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -0.99%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''139.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -0.91%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D28").Value = '''6.88'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +1.02%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''15.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +0.16%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''1.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  +0.07%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''0.0494'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -0.31%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''3.33'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +1.06%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''3.25'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +0.45%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = '''  +1.39%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = '''  -0.02%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''0.903'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -0.21%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = '''ImmutableX'
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = '''0.551'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -0.96%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = '''MXToken'
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = '''2.53'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -0.76%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''1.121.77'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -1.03%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '''  -0.35%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = '''  -0.12%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''5.59'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  +1.77%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''99.84'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +1.00%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  -0.39%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''1.779.46'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  +0.12%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.0₆0108'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -1.89%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''55.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -0.72%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''7.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -0.94%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''0.418'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -2.35%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = '''Cronos'
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = '''0.0503'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -0.06%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = '''SynthetixNetwork'
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = '''https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = '''2.36'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +3.13%  '
$ws.Range("E51").Style = "Normal"
